$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 13 (between Hakeem Olajuwon at row 12 and Joakim Noah at row 13)
$ws.Rows.Item(13).Insert()

$ws.Range("A13").Value = "Jaren Jackson Jr."
$ws.Range("B13").Value = 1
